$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ensure numeric-looking text values in column D stay as text (matches source data
# which stores prices as literal strings, e.g. "0.639", not as floating point numbers).
$ws.Range("D2").Value = '61.927.83'
$ws.Range("E2").Value = '  -0.66%  '
$ws.Range("D3").Value = '3.416.22'
$ws.Range("E3").Value = '  -0.74%  '
$ws.Range("E4").Value = '  -0.16%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '409.88'
$ws.Range("E5").Value = '  +0.59%  '
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = '129.03'
$ws.Range("E6").Value = '  -3.48%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.639'
$ws.Range("E7").Value = '  +7.87%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.999'
$ws.Range("E8").Value = '  -0.08%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.736'
$ws.Range("E9").Value = '  +7.26%  '
$ws.Range("E10").Value = '  +11.04%  '
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '42.87'
$ws.Range("E11").Value = '  +2.00%  '
$ws.Range("B12").Value = 'Polkadot'
$ws.Range("C12").Value = 'https://coinranking.com/coin/25W7FG7om+polkadot-dot'
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = '9.12'
$ws.Range("E12").Value = '  +7.55%  '
$ws.Range("B13").Value = 'TRON'
$ws.Range("C13").Value = 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx'
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '0.141'
$ws.Range("E13").Value = '  -0.31%  '
$ws.Range("D14").Value = '3.943.58'
$ws.Range("E14").Value = '  -0.91%  '
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = '21.19'
$ws.Range("E15").Value = '  +6.35%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '0.0000203'
$ws.Range("E16").Value = '  +41.82%  '
$ws.Range("D17").Value = '3.389.53'
$ws.Range("E17").Value = '  -1.06%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '12.13'
$ws.Range("E18").Value = '  +6.84%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '1.08'
$ws.Range("E19").Value = '  +5.76%  '
$ws.Range("D20").Value = '61.739.44'
$ws.Range("E20").Value = '  -1.06%  '
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = '457.26'
$ws.Range("E21").Value = '  +45.52%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '92.18'
$ws.Range("E22").Value = '  +10.14%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '3.16'
$ws.Range("E23").Value = '  -0.54%  '
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = '12.93'
$ws.Range("E24").Value = '  +0.57%  '
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = '3.25'
$ws.Range("E25").Value = '  +2.64%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '33.80'
$ws.Range("E26").Value = '  +13.94%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '8.79'
$ws.Range("E27").Value = '  +8.16%  '
$ws.Range("E28").Value = '  -0.85%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '7.60'
$ws.Range("E29").Value = '  -0.21%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '2.75'
$ws.Range("E30").Value = '  -0.41%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '12.01'
$ws.Range("E31").Value = '  +5.56%  '
$ws.Range("B32").Value = 'Hedera'
$ws.Range("C32").Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.114'
$ws.Range("E32").Value = '  -1.20%  '
$ws.Range("B33").Value = 'Kaspa'
$ws.Range("C33").Value = 'https://coinranking.com/coin/V8GxkwWow+kaspa-kas'
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '0.168'
$ws.Range("E33").Value = '  -3.05%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '42.84'
$ws.Range("E34").Value = '  +0.67%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '1.00'
$ws.Range("E35").Value = '  -0.02%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.0497'
$ws.Range("E36").Value = '  +2.33%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '53.39'
$ws.Range("E37").Value = '  +3.48%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '0.998'
$ws.Range("E38").Value = '  -0.08%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '3.37'
$ws.Range("E39").Value = '  -0.91%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.135'
$ws.Range("E40").Value = '  +7.63%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '2.91'
$ws.Range("E41").Value = '  -0.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.315'
$ws.Range("E42").Value = '  -2.54%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '140.51'
$ws.Range("E43").Value = '  +1.46%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '4.19'
$ws.Range("E44").Value = '  +5.27%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '1.98'
$ws.Range("E45").Value = '  -0.28%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.41'
$ws.Range("E46").Value = '  +8.20%  '
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '16.60'
$ws.Range("E47").Value = '  -1.19%  '
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '22.63'
$ws.Range("E48").Value = '  +5.73%  '
$ws.Range("D49").Value = '3.750.86'
$ws.Range("E49").Value = '  -0.88%  '
$ws.Range("D50").Value = '2.108.48'
$ws.Range("E50").Value = '  -0.68%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '107.13'
$ws.Range("E51").Value = '  +28.21%  '
